$d = $word.ActiveDocument

$replacements = @(
    @("78÷3=", "42÷7="),
    @("48÷6=", "22÷8="),
    @("23÷2=", "37÷3="),
    @("24÷9=", "51÷9="),
    @("37÷7=", "21÷6="),
    @("57÷6=", "75÷6="),
    @("78÷6=", "90÷3="),
    @("35÷7=", "97÷7="),
    @("64÷7=", "89÷5="),
    @("72÷8=", "22÷6="),
    @("62÷3=", "10÷4="),
    @("45÷7=", "68÷7="),
    @("17÷9=", "56÷2="),
    @("34÷2=", "93÷2="),
    @("24÷7=", "84÷4="),
    @("97÷3=", "30÷4="),
    @("52÷4=", "82÷8="),
    @("86÷7=", "50÷6="),
    @("95÷9=", "92÷7="),
    @("67÷7=", "62÷2="),
    @("17÷3=", "58÷5="),
    @("38÷4=", "65÷2="),
    @("21÷2=", "41÷7="),
    @("73÷5=", "39÷7="),
    @("88÷9=", "31÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
